# Fix model sheet -- 'type' now looks to match the prompt type and expands
# recursively. The generic JS-ish type names (string/number/object) plus the
# now-redundant "elementType" column are replaced by the real prompt type
# values (text/decimal/geopoint/select_one), matching the 'survey' sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# name -> type (was: name -> type -> elementType)
$ws.Range("B2").Value = "text"        # refrigerator_id        : string -> text
$ws.Range("B3").Value = "decimal"     # refrigerator_size       : number -> decimal
$ws.Range("B4").Value = "geopoint"    # refrigerator_location   : object/elementType=geopoint -> geopoint
$ws.Range("B5").Value = "select_one"  # refrigerator_condition  : string -> select_one
# refrigerator_stock_level / integer (row 6) is unchanged.

# Drop the now-unused elementType column (C1 header + C4 value).
$ws.Range("C1:C6").ClearContents()

# Match the author's final selection on the model sheet.
[void]$ws.Activate()
$ws.Range("B4").Select() | Out-Null
